$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: refresh an existing enquiry's reply date ---
$ws.Range("A2").Value = 1
$ws.Range("C2").Value = 3
$ws.Range("G2").Value = 45769.88321890046
$ws.Range("G2").NumberFormat = "yyyy-MM-dd HH:mm:ss"

# --- Row 12: new enquiry from T1234567J ---
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "T1234567J"
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = "Test"
$ws.Range("F12").Value = 45769.78912310185
$ws.Range("F12").NumberFormat = "yyyy-MM-dd HH:mm:ss"

# --- Row 13: new enquiry from T2345678D ---
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "T2345678D"
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = "test"
$ws.Range("F13").Value = 45769.89216858796
$ws.Range("F13").NumberFormat = "yyyy-MM-dd HH:mm:ss"
